$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 23810754
$ws.Range("I6").Value = 29411992
$ws.Range("K6").Value = 88235976
$ws.Range("M6").Value = -88235864
$ws.Range("H32").Value = 5530.1875
$ws.Range("I32").Value = 5378
$ws.Range("K32").Value = 5378
$ws.Range("M32").Value = -5052
$ws.Range("H38").Value = 62500120
$ws.Range("I38").Value = 138
$ws.Range("J38").Value = 500000000
$ws.Range("K38").Value = 414
$ws.Range("L38").Value = 1500000000
$ws.Range("M38").Value = -42
$ws.Range("N38").Value = -1500000744
$ws.Range("H40").Value = 2383.3333
$ws.Range("I40").Value = 2077.7778
$ws.Range("J40").Value = 2688.889
$ws.Range("K40").Value = 2077.7778
$ws.Range("L40").Value = 2688.889
$ws.Range("M40").Value = -1902.7778
$ws.Range("N40").Value = -3038.889
$ws.Range("H41").Value = 1100404.2
$ws.Range("J41").Value = 500195
$ws.Range("L41").Value = 500195
$ws.Range("N41").Value = -501075
$ws.Range("H44").Value = 27758.438
$ws.Range("I44").Value = 36033.75
$ws.Range("J44").Value = 25000
$ws.Range("K44").Value = 36033.75
$ws.Range("L44").Value = 25000
$ws.Range("M44").Value = -35571.75
$ws.Range("N44").Value = -25924
$ws.Range("H51").Value = 5649.1665
$ws.Range("I51").Value = 3200.25
$ws.Range("J51").Value = 6873.625
$ws.Range("K51").Value = 3200.25
$ws.Range("L51").Value = 6873.625
$ws.Range("M51").Value = -2716.25
$ws.Range("N51").Value = -7841.625
$ws.Range("H52").Value = 630
$ws.Range("I52").Value = 630
$ws.Range("K52").Value = 1890
$ws.Range("M52").Value = -1730
$ws.Range("H58").Value = 125002880
$ws.Range("I58").Value = 515
$ws.Range("K58").Value = 1545
$ws.Range("M58").Value = -1395
$ws.Range("H62").Value = 2626100.8
$ws.Range("J62").Value = 500000
$ws.Range("L62").Value = 500000
$ws.Range("N62").Value = -501248
$ws.Range("H65").Value = 2626100.8
$ws.Range("J65").Value = 500000
$ws.Range("L65").Value = 2500000
$ws.Range("N65").Value = -2506240
$ws.Range("H69").Value = 7751.3
$ws.Range("I69").Value = 6006.5
$ws.Range("J69").Value = 8187.5
$ws.Range("K69").Value = 18019.5
$ws.Range("L69").Value = 24562.5
$ws.Range("N69").Value = -26310.5
$ws.Range("M69").Value = -17145.5
$ws.Range("H72").Value = 7751.3
$ws.Range("I72").Value = 6006.5
$ws.Range("J72").Value = 8187.5
$ws.Range("K72").Value = 54058.5
$ws.Range("L72").Value = 73687.5
$ws.Range("N72").Value = -82423.5
$ws.Range("M72").Value = -49690.5
$ws.Range("H86").Value = 16755870
$ws.Range("I86").Value = 2861
$ws.Range("J86").Value = 25132376
$ws.Range("K86").Value = 2861
$ws.Range("L86").Value = 25132376
$ws.Range("M86").Value = -1738
$ws.Range("N86").Value = -25134622
$ws.Range("H89").Value = 16755870
$ws.Range("I89").Value = 2861
$ws.Range("J89").Value = 25132376
$ws.Range("K89").Value = 14305
$ws.Range("L89").Value = 125661880
$ws.Range("M89").Value = -8689
$ws.Range("N89").Value = -125673112
$ws.Range("H98").Value = 1354.0968
$ws.Range("I98").Value = 1097.1923
$ws.Range("K98").Value = 1097.1923
$ws.Range("M98").Value = 400.8077000000001
$ws.Range("H122").Value = 1354.0968
$ws.Range("I122").Value = 1097.1923
$ws.Range("K122").Value = 3291.5769
$ws.Range("M122").Value = -841.5769
$ws.Range("H125").Value = 1061.7858
$ws.Range("J125").Value = 1093.4
$ws.Range("L125").Value = 9840.6
$ws.Range("N125").Value = -14760.6
$ws.Range("H127").Value = 1617.8077
$ws.Range("I127").Value = 1316.4736
$ws.Range("K127").Value = 3949.4208
$ws.Range("M127").Value = 1010.5792
$ws.Range("H129").Value = 1594.7778
$ws.Range("I129").Value = 1310.1666
$ws.Range("J129").Value = 2164
$ws.Range("K129").Value = 3930.4998
$ws.Range("L129").Value = 6492
$ws.Range("M129").Value = 1069.5002
$ws.Range("N129").Value = -16492
$ws.Range("H132").Value = 4671.913
$ws.Range("I132").Value = 4271.0513
$ws.Range("K132").Value = 12813.1539
$ws.Range("M132").Value = -10283.1539
$ws.Range("H137").Value = 1999.4642
$ws.Range("I137").Value = 1715.8975
$ws.Range("K137").Value = 5147.6925
$ws.Range("M137").Value = -2597.6925
$ws.Range("H138").Value = 2886.3125
$ws.Range("I138").Value = 1174.1034
$ws.Range("J138").Value = 4305
$ws.Range("K138").Value = 3522.3102
$ws.Range("L138").Value = 12915
$ws.Range("M138").Value = 1617.6898
$ws.Range("N138").Value = -23195
$ws.Range("H141").Value = 2167.7727
$ws.Range("I141").Value = 1210.9445
$ws.Range("J141").Value = 6473.5
$ws.Range("K141").Value = 3632.8335
$ws.Range("L141").Value = 19420.5
$ws.Range("M141").Value = 1547.1665
$ws.Range("N141").Value = -29780.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1667.421
$ws.Range("I2").Value = 1763.0588
$ws.Range("J2").Value = 854.5
$ws.Range("K2").Value = 1763.0588
$ws.Range("L2").Value = 854.5
$ws.Range("M2").Value = -1650.0588
$ws.Range("N2").Value = -1080.5
$ws.Range("H32").Value = 192654.39
$ws.Range("I32").Value = 211665.83
$ws.Range("J32").Value = 10144.6
$ws.Range("K32").Value = 211665.83
$ws.Range("L32").Value = 10144.6
$ws.Range("M32").Value = -211378.83
$ws.Range("N32").Value = -10718.6
$ws.Range("H37").Value = 43333.332
$ws.Range("I37").Value = 15000
$ws.Range("K37").Value = 15000
$ws.Range("M37").Value = -14727
$ws.Range("H61").Value = 14494188
$ws.Range("I61").Value = 16130443
$ws.Range("J61").Value = 1645.7142
$ws.Range("K61").Value = 16130443
$ws.Range("L61").Value = 1645.7142
$ws.Range("M61").Value = -16130231
$ws.Range("N61").Value = -2069.7142
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H74").Value = 2846.4348
$ws.Range("I74").Value = 2696.2092
$ws.Range("J74").Value = 4999.6665
$ws.Range("K74").Value = 2696.2092
$ws.Range("L74").Value = 4999.6665
$ws.Range("M74").Value = -1822.2092
$ws.Range("N74").Value = -6747.6665
$ws.Range("H76").Value = 43429.332
$ws.Range("J76").Value = 43429.332
$ws.Range("L76").Value = 43429.332
$ws.Range("N76").Value = -44105.332
$ws.Range("H77").Value = 2846.4348
$ws.Range("I77").Value = 2696.2092
$ws.Range("J77").Value = 4999.6665
$ws.Range("K77").Value = 13481.046
$ws.Range("L77").Value = 24998.3325
$ws.Range("M77").Value = -9113.045999999998
$ws.Range("N77").Value = -33734.3325
$ws.Range("H79").Value = 43429.332
$ws.Range("J79").Value = 43429.332
$ws.Range("L79").Value = 43429.332
$ws.Range("N79").Value = -45769.332
$ws.Range("H88").Value = 8773118
$ws.Range("I88").Value = 33333818
$ws.Range("J88").Value = 1439.5714
$ws.Range("K88").Value = 33333818
$ws.Range("L88").Value = 1439.5714
$ws.Range("M88").Value = -33333412
$ws.Range("N88").Value = -2251.5714
$ws.Range("H91").Value = 8773118
$ws.Range("I91").Value = 33333818
$ws.Range("J91").Value = 1439.5714
$ws.Range("K91").Value = 33333818
$ws.Range("L91").Value = 1439.5714
$ws.Range("M91").Value = -33332414
$ws.Range("N91").Value = -4247.5714
$ws.Range("H110").Value = 6749.6113
$ws.Range("I110").Value = 7100.2
$ws.Range("K110").Value = 7100.2
$ws.Range("M110").Value = -5055.2
$ws.Range("H116").Value = 1667.421
$ws.Range("I116").Value = 1763.0588
$ws.Range("J116").Value = 854.5
$ws.Range("K116").Value = 1763.0588
$ws.Range("L116").Value = 854.5
$ws.Range("M116").Value = 530.9412
$ws.Range("N116").Value = -5442.5
$ws.Range("H132").Value = 2443.5806
$ws.Range("I132").Value = 2384.1428
$ws.Range("K132").Value = 7152.428400000001
$ws.Range("M132").Value = -4622.428400000001
$ws.Range("H136").Value = 14494188
$ws.Range("I136").Value = 16130443
$ws.Range("J136").Value = 1645.7142
$ws.Range("K136").Value = 48391329
$ws.Range("L136").Value = 4937.142599999999
$ws.Range("M136").Value = -48388779
$ws.Range("N136").Value = -10037.1426
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1667.421
$ws.Range("I3").Value = 1763.0588
$ws.Range("J3").Value = 854.5
$ws.Range("K3").Value = 1763.0588
$ws.Range("L3").Value = 854.5
$ws.Range("M3").Value = -1649.0588
$ws.Range("N3").Value = -1082.5
$ws.Range("H20").Value = 17288.818
$ws.Range("I20").Value = 16310.706
$ws.Range("J20").Value = 20614.4
$ws.Range("K20").Value = 16310.706
$ws.Range("L20").Value = 20614.4
$ws.Range("M20").Value = -16063.706
$ws.Range("N20").Value = -21108.4
$ws.Range("H86").Value = 3917.7058
$ws.Range("I86").Value = 3946.7334
$ws.Range("K86").Value = 3946.7334
$ws.Range("M86").Value = -2823.7334
$ws.Range("H89").Value = 3917.7058
$ws.Range("I89").Value = 3946.7334
$ws.Range("K89").Value = 19733.667
$ws.Range("M89").Value = -14117.667
$ws.Range("H94").Value = 3103.625
$ws.Range("I94").Value = 3434.8333
$ws.Range("K94").Value = 3434.8333
$ws.Range("M94").Value = -2983.8333
$ws.Range("H103").Value = 39299.8
$ws.Range("J103").Value = 39299.8
$ws.Range("L103").Value = 39299.8
$ws.Range("N103").Value = -41643.8
$ws.Range("H105").Value = 2785.0476
$ws.Range("I105").Value = 3499
$ws.Range("K105").Value = 3499
$ws.Range("M105").Value = -1752
$ws.Range("H107").Value = 901
$ws.Range("I107").Value = 901
$ws.Range("K107").Value = 901
$ws.Range("M107").Value = 1019
$ws.Range("H134").Value = 1778.6364
$ws.Range("I134").Value = 1438.0667
$ws.Range("J134").Value = 2508.4285
$ws.Range("K134").Value = 4314.2001
$ws.Range("L134").Value = 7525.2855
$ws.Range("M134").Value = -1779.2001
$ws.Range("N134").Value = -12595.2855
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2166826.2
$ws.Range("J4").Value = 2166826.2
$ws.Range("L4").Value = 2166826.2
$ws.Range("N4").Value = -2167050.2
$ws.Range("H16").Value = 2264.2
$ws.Range("J16").Value = 2332.6667
$ws.Range("L16").Value = 2332.6667
$ws.Range("N16").Value = -2906.6667
$ws.Range("H22").Value = 348.875
$ws.Range("I22").Value = 358.73077
$ws.Range("J22").Value = 330.57144
$ws.Range("K22").Value = 358.73077
$ws.Range("L22").Value = 330.57144
$ws.Range("M22").Value = -8.730770000000007
$ws.Range("N22").Value = -1030.57144
$ws.Range("H31").Value = 4985.136
$ws.Range("I31").Value = 16788.25
$ws.Range("J31").Value = 3357.1206
$ws.Range("K31").Value = 16788.25
$ws.Range("L31").Value = 3357.1206
$ws.Range("M31").Value = -16493.25
$ws.Range("N31").Value = -3947.1206
$ws.Range("H34").Value = 4985.136
$ws.Range("I34").Value = 16788.25
$ws.Range("J34").Value = 3357.1206
$ws.Range("K34").Value = 16788.25
$ws.Range("L34").Value = 3357.1206
$ws.Range("M34").Value = -16586.25
$ws.Range("N34").Value = -3761.1206
$ws.Range("H58").Value = 502164.06
$ws.Range("I58").Value = 1383.4546
$ws.Range("J58").Value = 1114229.2
$ws.Range("K58").Value = 1383.4546
$ws.Range("L58").Value = 1114229.2
$ws.Range("M58").Value = -1180.4546
$ws.Range("N58").Value = -1114635.2
$ws.Range("H62").Value = 6606.643
$ws.Range("I62").Value = 6345.615
$ws.Range("K62").Value = 6345.615
$ws.Range("M62").Value = -5721.615
$ws.Range("H65").Value = 6606.643
$ws.Range("I65").Value = 6345.615
$ws.Range("K65").Value = 31728.075
$ws.Range("M65").Value = -28608.075
$ws.Range("H105").Value = 2653.7334
$ws.Range("I105").Value = 1507
$ws.Range("J105").Value = 3964.2856
$ws.Range("K105").Value = 1507
$ws.Range("L105").Value = 3964.2856
$ws.Range("M105").Value = 240
$ws.Range("N105").Value = -7458.2856
$ws.Range("H107").Value = 7241
$ws.Range("I107").Value = 5999.8335
$ws.Range("K107").Value = 5999.8335
$ws.Range("M107").Value = -4079.8335
$ws.Range("H113").Value = 2264.2
$ws.Range("J113").Value = 2332.6667
$ws.Range("L113").Value = 2332.6667
$ws.Range("N113").Value = -6672.6667
$ws.Range("H134").Value = 1090.4584
$ws.Range("I134").Value = 884.6316
$ws.Range("K134").Value = 2653.8948
$ws.Range("M134").Value = -118.8948
$ws.Range("H136").Value = 502164.06
$ws.Range("I136").Value = 1383.4546
$ws.Range("J136").Value = 1114229.2
$ws.Range("K136").Value = 4150.3638
$ws.Range("L136").Value = 3342687.6
$ws.Range("M136").Value = -1600.3638
$ws.Range("N136").Value = -3347787.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 65.42856999999999
$ws.Range("I2").Value = 42.05263
$ws.Range("J2").Value = 84.73913
$ws.Range("K2").Value = 252.31578
$ws.Range("L2").Value = 508.43478
$ws.Range("M2").Value = -139.31578
$ws.Range("N2").Value = -734.43478
$ws.Range("H12").Value = 742.5789
$ws.Range("I12").Value = 497.2857
$ws.Range("J12").Value = 885.6667
$ws.Range("K12").Value = 1491.8571
$ws.Range("L12").Value = 2657.0001
$ws.Range("M12").Value = -1318.8571
$ws.Range("N12").Value = -3003.0001
$ws.Range("H22").Value = 9749.375
$ws.Range("I22").Value = 8999
$ws.Range("K22").Value = 26997
$ws.Range("M22").Value = -26828
$ws.Range("H27").Value = 9749.375
$ws.Range("I27").Value = 8999
$ws.Range("K27").Value = 26997
$ws.Range("M27").Value = -26895
$ws.Range("H50").Value = 116.8
$ws.Range("I50").Value = 116.8
$ws.Range("K50").Value = 350.4
$ws.Range("M50").Value = 130.6
$ws.Range("H53").Value = 116.8
$ws.Range("I53").Value = 116.8
$ws.Range("K53").Value = 350.4
$ws.Range("M53").Value = 130.6
$ws.Range("H64").Value = 175644.14
$ws.Range("I64").Value = 243802
$ws.Range("J64").Value = 5249.5
$ws.Range("K64").Value = 731406
$ws.Range("L64").Value = 15748.5
$ws.Range("M64").Value = -731136
$ws.Range("N64").Value = -16288.5
$ws.Range("H67").Value = 175644.14
$ws.Range("I67").Value = 243802
$ws.Range("J67").Value = 5249.5
$ws.Range("K67").Value = 731406
$ws.Range("L67").Value = 15748.5
$ws.Range("M67").Value = -730470
$ws.Range("N67").Value = -17620.5
$ws.Range("H80").Value = 5972.577
$ws.Range("I80").Value = 5635.1763
$ws.Range("J80").Value = 6609.8887
$ws.Range("K80").Value = 16905.5289
$ws.Range("L80").Value = 19829.6661
$ws.Range("M80").Value = -15969.5289
$ws.Range("N80").Value = -21701.6661
$ws.Range("H82").Value = 15014.595
$ws.Range("I82").Value = 15000
$ws.Range("K82").Value = 45000
$ws.Range("M82").Value = -44594
$ws.Range("H83").Value = 5972.577
$ws.Range("I83").Value = 5635.1763
$ws.Range("J83").Value = 6609.8887
$ws.Range("K83").Value = 50716.5867
$ws.Range("L83").Value = 59488.99830000001
$ws.Range("M83").Value = -46036.5867
$ws.Range("N83").Value = -68848.99830000001
$ws.Range("H85").Value = 15014.595
$ws.Range("I85").Value = 15000
$ws.Range("K85").Value = 45000
$ws.Range("M85").Value = -43596
$ws.Range("H99").Value = 1405.75
$ws.Range("I99").Value = 1405.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4217.25
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1971.25
$ws.Range("N99").ClearContents()
$ws.Range("H103").Value = 519
$ws.Range("J103").Value = 1176
$ws.Range("L103").Value = 3528
$ws.Range("N103").Value = -5286
$ws.Range("H113").Value = 1093.2778
$ws.Range("I113").Value = 450
$ws.Range("J113").Value = 1221.9333
$ws.Range("K113").Value = 1350
$ws.Range("L113").Value = 3665.7999
$ws.Range("M113").Value = 820
$ws.Range("N113").Value = -8005.7999
$ws.Range("H117").Value = 710
$ws.Range("J117").Value = 840
$ws.Range("L117").Value = 2520
$ws.Range("N117").Value = -9404
$ws.Range("H121").Value = 481898
$ws.Range("I121").Value = 674027.5600000001
$ws.Range("J121").Value = 1574.1666
$ws.Range("K121").Value = 2022082.68
$ws.Range("L121").Value = 4722.4998
$ws.Range("M121").Value = -2020772.68
$ws.Range("N121").Value = -7342.4998
$ws.Range("H139").Value = 1968610.8
$ws.Range("I139").Value = 2091492.6
$ws.Range("J139").Value = 2500
$ws.Range("K139").Value = 6274477.800000001
$ws.Range("L139").Value = 7500
$ws.Range("M139").Value = -6269337.800000001
$ws.Range("N139").Value = -17780
$ws.Range("H140").Value = 1572.875
$ws.Range("I140").Value = 1321.4286
$ws.Range("K140").Value = 3964.2858
$ws.Range("M140").Value = 1215.7142
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 2501900
$ws.Range("I29").Value = 2501900
$ws.Range("K29").Value = 2501900
$ws.Range("M29").Value = -2501610
$ws.Range("H35").Value = 26000
$ws.Range("J35").Value = 26000
$ws.Range("L35").Value = 26000
$ws.Range("N35").Value = -26596
$ws.Range("H70").Value = 14863.333
$ws.Range("I70").Value = 14825
$ws.Range("K70").Value = 14825
$ws.Range("M70").Value = -14555
$ws.Range("H73").Value = 14863.333
$ws.Range("I73").Value = 14825
$ws.Range("K73").Value = 14825
$ws.Range("M73").Value = -13889
$ws.Range("H97").Value = 1215.4
$ws.Range("I97").Value = 1280.375
$ws.Range("J97").Value = 955.5
$ws.Range("K97").Value = 1280.375
$ws.Range("L97").Value = 955.5
$ws.Range("M97").Value = -784.375
$ws.Range("N97").Value = -1947.5
$ws.Range("H107").Value = 2331.35
$ws.Range("I107").Value = 2591.1428
$ws.Range("K107").Value = 2591.1428
$ws.Range("M107").Value = -671.1428000000001
$ws.Range("H113").Value = 6992.5
$ws.Range("I113").Value = 6425.857
$ws.Range("J113").Value = 8314.666999999999
$ws.Range("K113").Value = 6425.857
$ws.Range("L113").Value = 8314.666999999999
$ws.Range("M113").Value = -4255.857
$ws.Range("N113").Value = -12654.667
$ws.Range("H132").Value = 406094.38
$ws.Range("I132").Value = 561296.4
$ws.Range("J132").Value = 7003.4287
$ws.Range("K132").Value = 1683889.2
$ws.Range("L132").Value = 21010.2861
$ws.Range("M132").Value = -1681359.2
$ws.Range("N132").Value = -26070.2861
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1833.7222
$ws.Range("I7").Value = 1833.7222
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1833.7222
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1721.7222
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 936.25
$ws.Range("I22").Value = 1003.1111
$ws.Range("K22").Value = 1003.1111
$ws.Range("M22").Value = -708.1111
$ws.Range("H27").Value = 936.25
$ws.Range("I27").Value = 1003.1111
$ws.Range("K27").Value = 1003.1111
$ws.Range("M27").Value = -896.1111
$ws.Range("H32").Value = 2946.2856
$ws.Range("I32").Value = 2946.2856
$ws.Range("K32").Value = 2946.2856
$ws.Range("M32").Value = -2629.2856
$ws.Range("H46").Value = 1844.0555
$ws.Range("I46").Value = 1699.7273
$ws.Range("K46").Value = 1699.7273
$ws.Range("M46").Value = -1511.7273
$ws.Range("H93").Value = 6562.25
$ws.Range("I93").Value = 4000
$ws.Range("J93").Value = 8099.6
$ws.Range("K93").Value = 4000
$ws.Range("L93").Value = 8099.6
$ws.Range("M93").Value = -2752
$ws.Range("N93").Value = -10595.6
$ws.Range("H100").Value = 1996.6666
$ws.Range("I100").Value = 1995
$ws.Range("K100").Value = 1995
$ws.Range("M100").Value = -1454
$ws.Range("H109").Value = 500001
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("H111").Value = 55600
$ws.Range("J111").Value = 55600
$ws.Range("L111").Value = 55600
$ws.Range("N111").Value = -63780
$ws.Range("H122").Value = 6263.4287
$ws.Range("I122").Value = 4198.6665
$ws.Range("K122").Value = 12595.9995
$ws.Range("M122").Value = -10145.9995
$ws.Range("H126").Value = 1833.7222
$ws.Range("I126").Value = 1833.7222
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5501.1666
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3031.1666
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3452.5208
$ws.Range("I132").Value = 1497.4546
$ws.Range("J132").Value = 7753.6665
$ws.Range("K132").Value = 4492.3638
$ws.Range("L132").Value = 23260.9995
$ws.Range("M132").Value = -1962.3638
$ws.Range("N132").Value = -28320.9995
$ws.Range("H136").Value = 2112.5193
$ws.Range("I136").Value = 1789.8334
$ws.Range("J136").Value = 3467.8
$ws.Range("K136").Value = 5369.5002
$ws.Range("L136").Value = 10403.4
$ws.Range("M136").Value = -2819.5002
$ws.Range("N136").Value = -15503.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 35005
$ws.Range("J24").Value = 35005
$ws.Range("L24").Value = 35005
$ws.Range("N24").Value = -35465
$ws.Range("H62").Value = 1237975.9
$ws.Range("J62").Value = 67614.766
$ws.Range("L62").Value = 67614.766
$ws.Range("N62").Value = -68862.766
$ws.Range("H65").Value = 1237975.9
$ws.Range("J65").Value = 67614.766
$ws.Range("L65").Value = 338073.83
$ws.Range("N65").Value = -344313.83
$ws.Range("H68").Value = 43045.285
$ws.Range("J68").Value = 43045.285
$ws.Range("L68").Value = 43045.285
$ws.Range("N68").Value = -44667.285
$ws.Range("H71").Value = 43045.285
$ws.Range("J71").Value = 43045.285
$ws.Range("L71").Value = 129135.855
$ws.Range("N71").Value = -137247.855
$ws.Range("H82").Value = 75000
$ws.Range("J82").Value = 75000
$ws.Range("L82").Value = 75000
$ws.Range("N82").Value = -75766
$ws.Range("H85").Value = 75000
$ws.Range("J85").Value = 75000
$ws.Range("L85").Value = 75000
$ws.Range("N85").Value = -77652
$ws.Range("H100").Value = 4801.5625
$ws.Range("J100").Value = 1448
$ws.Range("L100").Value = 2896
$ws.Range("N100").Value = -3978
$ws.Range("H107").Value = 1175.25
$ws.Range("I107").Value = 500.5
$ws.Range("J107").Value = 1850
$ws.Range("K107").Value = 1501.5
$ws.Range("L107").Value = 5550
$ws.Range("M107").Value = 418.5
$ws.Range("N107").Value = -9390
$ws.Range("H113").Value = 7093.421
$ws.Range("I113").Value = 8827.666999999999
$ws.Range("K113").Value = 26483.001
$ws.Range("M113").Value = -24313.001
$ws.Range("H119").Value = 53699.8
$ws.Range("J119").Value = 53999.75
$ws.Range("L119").Value = 53999.75
$ws.Range("N119").Value = -63675.75
$ws.Range("H126").Value = 1414.4445
$ws.Range("I126").Value = 1144
$ws.Range("K126").Value = 3432
$ws.Range("M126").Value = -962
$ws.Range("H132").Value = 576314
$ws.Range("I132").Value = 773161.5600000001
$ws.Range("J132").Value = 7643.3335
$ws.Range("K132").Value = 2319484.68
$ws.Range("L132").Value = 22930.0005
$ws.Range("M132").Value = -2316954.68
$ws.Range("N132").Value = -27990.0005
